$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.456.89'
$ws.Range("E2").Value = '  +1.19%  '

# Row 3
$ws.Range("D3").Value = '1.825.95'
$ws.Range("E3").Value = '  +1.65%  '

# Row 4
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").Value = '316.27'
$ws.Range("E5").Value = '  -0.25%  '

# Row 6
$ws.Range("D6").Value = '1.000'

# Row 7
$ws.Range("D7").Value = '0.5431'
$ws.Range("E7").Value = '  +1.59%  '

# Row 8
$ws.Range("D8").Value = '0.4033'
$ws.Range("E8").Value = '  +7.01%  '

# Row 9
$ws.Range("D9").Value = '0.07666'
$ws.Range("E9").Value = '  +2.76%  '

# Row 10
$ws.Range("D10").Value = '41.87'
$ws.Range("E10").Value = '  +0.13%  '

# Row 11
$ws.Range("E11").Value = '  +1.45%  '

# Row 12
$ws.Range("D12").Value = '6.327'
$ws.Range("E12").Value = '  +3.45%  '

# Row 13
$ws.Range("D13").Value = '7.653'
$ws.Range("E13").Value = '  +5.72%  '

# Row 14
$ws.Range("D14").Value = '1.000'
$ws.Range("E14").Value = '  +0.02%  '

# Row 15
$ws.Range("D15").Value = '20.93'
$ws.Range("E15").Value = '  +1.23%  '

# Row 16
$ws.Range("D16").Value = '1.828.64'
$ws.Range("E16").Value = '  +1.96%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001078'
$ws.Range("E17").Value = '  +2.08%  '

# Row 18
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '89.68'
$ws.Range("E18").Value = '  +0.64%  '

# Row 19
$ws.Range("D19").Value = '0.06595'
$ws.Range("E19").Value = '  +2.05%  '

# Row 20
$ws.Range("D20").Value = '17.65'
$ws.Range("E20").Value = '  +1.56%  '

# Row 21
$ws.Range("E21").Value = '  +0.24%  '

# Row 22
$ws.Range("D22").Value = '6.069'
$ws.Range("E22").Value = '  +2.71%  '

# Row 23
$ws.Range("D23").Value = '28.462.09'
$ws.Range("E23").Value = '  +1.21%  '

# Row 24
$ws.Range("D24").Value = '11.09'
$ws.Range("E24").Value = '  -0.90%  '

# Row 25
$ws.Range("D25").Value = '2.228'
$ws.Range("E25").Value = '  +6.05%  '

# Row 26
$ws.Range("D26").Value = '2.467'
$ws.Range("E26").Value = '  +7.31%  '

# Row 27
$ws.Range("D27").Value = '20.74'
$ws.Range("E27").Value = '  +2.32%  '

# Row 28
$ws.Range("D28").Value = '157.17'
$ws.Range("E28").Value = '  +1.31%  '

# Row 29
$ws.Range("D29").Value = '2.038.72'
$ws.Range("E29").Value = '  +2.03%  '

# Row 30
$ws.Range("D30").Value = '123.83'
$ws.Range("E30").Value = '  +2.79%  '

# Row 31
$ws.Range("D31").Value = '0.1116'
$ws.Range("E31").Value = '  +6.43%  '

# Row 32
$ws.Range("D32").Value = '1.126'
$ws.Range("E32").Value = '  +0.69%  '

# Row 33
$ws.Range("E33").Value = '  +1.95%  '

# Row 34
$ws.Range("D34").Value = '0.07374'
$ws.Range("E34").Value = '  +13.66%  '

# Row 35
$ws.Range("D35").Value = '3.640'
$ws.Range("E35").Value = '  -0.33%  '

# Row 36
$ws.Range("E36").Value = '  -0.64%  '

# Row 37
$ws.Range("E37").Value = '  +2.28%  '

# Row 38
$ws.Range("D38").Value = '5.202'
$ws.Range("E38").Value = '  +3.48%  '

# Row 39
$ws.Range("D39").Value = '8.824'
$ws.Range("E39").Value = '  +3.80%  '

# Row 40
$ws.Range("D40").Value = '11.35'
$ws.Range("E40").Value = '  +2.40%  '

# Row 41
$ws.Range("D41").Value = '0.6271'
$ws.Range("E41").Value = '  +1.69%  '

# Row 42
$ws.Range("D42").Value = '1.176'
$ws.Range("E42").Value = '  +0.14%  '

# Row 43
$ws.Range("E43").Value = '  +0.25%  '

# Row 44
$ws.Range("E44").Value = '  -3.88%  '

# Row 45
$ws.Range("D45").Value = '13.51'
$ws.Range("E45").Value = '  +1.17%  '

# Row 46
$ws.Range("E46").Value = '  +0.65%  '

# Row 47
$ws.Range("D47").Value = '0.5847'
$ws.Range("E47").Value = '  +1.30%  '

# Row 48
$ws.Range("D48").Value = '124.92'
$ws.Range("E48").Value = '  -1.91%  '

# Row 49
$ws.Range("E49").Value = '  +3.86%  '

# Row 50
$ws.Range("D50").Value = '1.201'
$ws.Range("E50").Value = '  +0.92%  '

# Row 51
$ws.Range("E51").Value = '  +1.08%  '
